$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.353672031788087, 0.1490953738153317, 0.2416789000905896, 0.2456626311167568, 0.7152945399284363, 0.7446634769439697, 0.7068678140640259, 0.7268773913383484)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
